# Qatar Stars League workbook update
# - Swap the B:AC data between certain row pairs (the "id" counter in column A
#   stays put, but the match record that was attached to it moves to the
#   sibling row and vice versa).
# - Append six new upcoming fixtures as rows 86-91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, $rowA, $rowB)
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")
    $valA = $rangeA.Value()
    $valB = $rangeB.Value()
    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

$pairs = @(
    @(15,16),
    @(21,22),
    @(27,28),
    @(29,30),
    @(42,43),
    @(56,57),
    @(60,61),
    @(62,63),
    @(70,71),
    @(83,84)
)

foreach ($p in $pairs) {
    Swap-RowData $ws $p[0] $p[1]
}

# Copy the formatting (font/border/number format) of the last existing data
# row down onto the six freshly appended rows so the new cells pick up the
# same styles (bold/boxed id column, date-formatted Date column, etc.).
$ws.Range("A85:AC85").Copy()
$ws.Range("A86:AC91").PasteSpecial(-4122)
$excel.CutCopyMode = 0

function Set-NewRow {
    param($ws, $r, $vals)
    foreach ($key in $vals.Keys) {
        $ws.Range("$key$r").Value = $vals[$key]
    }
}

Set-NewRow $ws 86 @{
    A = 84; B = 7840686; C = "Qatar Stars League"; D = "Qatar Stars League";
    E = 45356.5; F = "Al Duhail"; G = "Al Markhiya";
    K = 1.285; L = 5.5; M = 7.5; N = 1.285; O = 5.5; P = 7.5;
    Q = -1.75; R = 1.975; S = 1.825; T = 3.5; U = 1.85; V = 1.95;
    W = 0; X = 0; Y = 0; Z = 0; AA = 0
}

Set-NewRow $ws 87 @{
    A = 85; B = 7840807; C = "Qatar Stars League"; D = "Qatar Stars League";
    E = 45356.58333333334; F = "AlMuaidar"; G = "AlAhli Doha";
    K = 2.875; L = 4; M = 2; N = 2.875; O = 4; P = 2;
    Q = 0.25; R = 1.975; S = 1.825; T = 3.5; U = 1.85; V = 1.95;
    W = 0; X = 0; Y = 0; Z = 0; AA = 0
}

Set-NewRow $ws 88 @{
    A = 86; B = 7840806; C = "Qatar Stars League"; D = "Qatar Stars League";
    E = 45356.58333333334; F = "Umm Salal"; G = "AlShamal SC";
    K = 2; L = 3.6; M = 3.2; N = 2; O = 3.6; P = 3.2;
    Q = -0.25; R = 1.775; S = 2.025; T = 3; U = 2; V = 1.8;
    W = 0; X = 0; Y = 0; Z = 0; AA = 0
}

Set-NewRow $ws 89 @{
    A = 87; B = 7840808; C = "Qatar Stars League"; D = "Qatar Stars League";
    E = 45357.5; F = "Qatar SC Doha"; G = "AlArabi Doha";
    K = 3.6; L = 3.6; M = 1.85; N = 3.6; O = 3.6; P = 1.85;
    Q = 0.5; R = 1.9; S = 1.9; T = 2.75; U = 1.8; V = 2;
    W = 0; X = 0; Y = 0; Z = 0; AA = 0
}

Set-NewRow $ws 90 @{
    A = 88; B = 7840810; C = "Qatar Stars League"; D = "Qatar Stars League";
    E = 45357.58333333334; F = "AlRayyan SC"; G = "AlWakrah SC";
    K = 2; L = 3.6; M = 3.1; N = 2; O = 3.6; P = 3.1;
    Q = -0.5; R = 2.025; S = 1.775; T = 3.25; U = 1.975; V = 1.825;
    W = 0; X = 0; Y = 0; Z = 0; AA = 0
}

Set-NewRow $ws 91 @{
    A = 89; B = 7840809; C = "Qatar Stars League"; D = "Qatar Stars League";
    E = 45357.58333333334; F = "Al Gharafa"; G = "Al Sadd";
    K = 5; L = 4.75; M = 1.45; N = 4.5; O = 4.5; P = 1.5;
    Q = 1.25; R = 1.8; S = 2; T = 3.75; U = 1.9; V = 1.9;
    W = 0; X = 0; Y = 0; Z = 0; AA = 0
}

# The formats-only copy above also duplicated the (empty) FTHG/FTAG/FTR and
# PL_AhOver/PL_AhUnder cells from the template row onto the new fixtures.
# Those matches have not been played yet, so those columns must stay
# completely absent rather than present-but-empty.
$ws.Range("H86:J91").ClearContents()
$ws.Range("AB86:AC91").ClearContents()
